$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.233.48'
$ws.Range('E2').Value = '  -0.67%  '

$ws.Range('D3').Value = '1.839.73'
$ws.Range('E3').Value = '  -1.42%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4651'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.85%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2724'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.47%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06274'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.08%  '

$ws.Range('D10').Value = '1.834.38'
$ws.Range('E10').Value = '  -1.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07415'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.38%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.16'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.935'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '83.48'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E15').Value = '  -3.18%  '

$ws.Range('D16').Value = '30.164.76'
$ws.Range('E16').Value = '  -0.85%  '

$ws.Range('E17').Value = '  +0.02%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.59%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007287'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.26%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.07%  '

$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.883'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.70%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.847'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.84%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.178'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.48%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.864'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1032'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.372'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.071'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.38%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.804'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04818'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.05%  '

$ws.Range('E33').Value = '  -2.78%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7034'
$ws.Range('D34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.690'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.69%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01866'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.90%  '

$ws.Range('E37').Value = '  +0.46%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8874'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.83%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '104.69'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.37%  '

$ws.Range('E40').Value = '  -5.91%  '

$ws.Range('E41').Value = '  +0.59%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.515'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4010'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.027'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.82%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.15'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.62%  '

$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1194'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.36%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.605'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.15%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.46%  '

$ws.Range('E49').Value = '  -2.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.350'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.28%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3627'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.75%  '
